# Automatische test-sync: 2025-06-29 14:53:50
# Appends the 6th test-mail log entry to the "Logs" sheet and refreshes the
# "Overig" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Logs" sheet: append new row 21 with the testmail #6 data.
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 21

$logs.Range("A$newRow").Value = "Kun je dit intern overleggen?"
$logs.Range("B$newRow").Value = "mailmind.test@zohomail.eu"
$logs.Range("C$newRow").Value = "Testmail #6: Kun je dit intern overleggen?"
$logs.Range("D$newRow").Value = "Overig"
$logs.Range("E$newRow").Value = "Beste afzender,`r`nBedankt voor je e-mail. Kun je wat meer specifieke informatie geven over waarover je precies wilt dat er intern overlegd wordt? Op die manier kan ik ervoor zorgen dat je aanvraag bij de juiste persoon of afdeling terechtkomt.`r`nMet vriendelijke groet,`r`n[Naam]`r`nE-mailassistent"
$logs.Range("F$newRow").Value = "2025-06-29 14:53:12"
$logs.Range("G$newRow").Value = "Ja"
$logs.Range("H$newRow").Value = "Ja"
$logs.Range("I$newRow").Value = "Nee"

# Extend the conditional-formatting ranges (D, G, H, I) so the newly added
# row is covered as well, matching how Excel keeps them in sync with the
# growing table.
$colsToExtend = "D", "G", "H", "I"
foreach ($col in $colsToExtend) {
    $rng = $logs.Range("$col" + "2:" + "$col" + "20")
    $fcs = $rng.FormatConditions
    $count = $fcs.Count()
    for ($i = 1; $i -le $count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($logs.Range("$col" + "2:" + "$col" + "$newRow"))
    }
}

# ---------------------------------------------------------------------
# 2. "Dashboard" sheet: bump the "Overig" count from 1 to 2.
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B6").Value = 2
